# Acronym key workbook update:
#  - Add a new "BCoESC" (BAU Cost of Electricity Sector Capital) row to the
#    "Key to Variables" sheet, in its correct alphabetically-sorted spot
#    within the "elec" section (just before "BCpUC").
#  - Remove the now-superseded "CoESC" (Cost of Electricity Sector Capital)
#    row from the same section.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Key to Variables")

# --- Insert the new row before row 62 ("BCpUC") -----------------------
$ws.Rows("62:62").Insert()

# Copy formatting from a same-shaped existing row (A:D wrap-text cells +
# an F cell styled for "medium") so the new row matches the sheet's look.
# Row 89 is "BPTBfRN"'s old neighbour row, now shifted down to 89 by the
# insert above, and has exactly the A/B/C/D + F cell layout we need.
$ws.Range("A89:D89").Copy()
$ws.Range("A62:D62").PasteSpecial(-4122)
$ws.Range("F89").Copy()
$ws.Range("F62").PasteSpecial(-4122)

$ws.Range("A62").Value = "elec"
$ws.Range("B62").Value = "BCoESC"
$ws.Range("C62").Value = "BAU Cost of Electricity Sector Capital"
$ws.Range("D62").Value = "BAU Cost of Electricity Sector Capital for Power Plants, BAU Cost of Electricity Sector Capital for Other Investments"
$ws.Range("F62").Value = "medium"
$ws.Rows("62:62").RowHeight = 30

# --- Remove the old "CoESC" row ----------------------------------------
$foundRow = -1
for ($r = 63; $r -le 90; $r++) {
    $bval = $ws.Cells.Item($r, 2).Value2
    if ($bval -eq "CoESC") {
        $foundRow = $r
    }
}
if ($foundRow -gt 0) {
    $rowRef = "$foundRow`:$foundRow"
    $ws.Rows($rowRef).Delete()
}

# --- Restore the view to roughly where the edit happened ---------------
$ws.Range("D63").Select()

# The workbook was (and remains) saved with the "About" sheet as the
# active/visible tab - re-activate it so we don't leave "Key to Variables"
# selected as a side-effect of editing it above.
$wb.Worksheets.Item("About").Activate()
